# Apply updated res_line/pl_mw.xlsx results for the 380 kV case.
# Updates columns B, D, E, F, G, I, J, K, O for data rows 2-25
# (columns C, H, L, M, N remain 0; column A is the row index).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1032961280896671
$ws.Range("D2").Value = 0.1973786334130665
$ws.Range("E2").Value = 0.1686580250957377
$ws.Range("F2").Value = 1.38539536047562
$ws.Range("G2").Value = 0.002456613967747282
$ws.Range("I2").Value = 0.4947891773267941
$ws.Range("J2").Value = 0.186744779948846
$ws.Range("K2").Value = 0.9426148357079001
$ws.Range("O2").Value = 3.30716897825846

# Row 3
$ws.Range("B3").Value = 0.09077806294087054
$ws.Range("D3").Value = 0.1913542004115669
$ws.Range("E3").Value = 0.1644126951552778
$ws.Range("F3").Value = 1.388371444129064
$ws.Range("G3").Value = 0.002459499932330808
$ws.Range("I3").Value = 0.5033369029767574
$ws.Range("J3").Value = 0.1827546154448214
$ws.Range("K3").Value = 0.8361995459853517
$ws.Range("O3").Value = 3.328628647878759

# Row 4
$ws.Range("B4").Value = 0.08307687039801692
$ws.Range("D4").Value = 0.1877237933466631
$ws.Range("E4").Value = 0.1618874637718051
$ws.Range("F4").Value = 1.391069844476128
$ws.Range("G4").Value = 0.002461366423824715
$ws.Range("I4").Value = 0.5088964811680303
$ws.Range("J4").Value = 0.1804138514620774
$ws.Range("K4").Value = 0.7707506856914961
$ws.Range("O4").Value = 3.344019255003417

# Row 5
$ws.Range("B5").Value = 0.07993502248385198
$ws.Range("D5").Value = 0.1862617322926781
$ws.Range("E5").Value = 0.1608789242480633
$ws.Range("F5").Value = 1.392388369240187
$ws.Range("G5").Value = 0.002462150868321077
$ws.Range("I5").Value = 0.5112403007616479
$ws.Range("J5").Value = 0.1794874291619308
$ws.Range("K5").Value = 0.744053757952031
$ws.Range("O5").Value = 3.350847233796259

# Row 6
$ws.Range("B6").Value = 0.07941311360033865
$ws.Range("D6").Value = 0.1860200097129194
$ws.Range("E6").Value = 0.1607126971189921
$ws.Range("F6").Value = 1.392620526686933
$ws.Range("G6").Value = 0.002462282566307344
$ws.Range("I6").Value = 0.5116342156060005
$ws.Range("J6").Value = 0.1793352557213339
$ws.Range("K6").Value = 0.7396192293490174
$ws.Range("O6").Value = 3.352014588387888

# Row 7
$ws.Range("B7").Value = 0.08303451239703463
$ws.Range("D7").Value = 0.1877040050533338
$ws.Range("E7").Value = 0.1618737791467275
$ws.Range("F7").Value = 1.391086740405214
$ws.Range("G7").Value = 0.002461376906675278
$ws.Range("I7").Value = 0.5089277739450209
$ws.Range("J7").Value = 0.18040124623586
$ws.Range("K7").Value = 0.7703907447619542
$ws.Range("O7").Value = 3.344109088383902

# Row 8
$ws.Range("B8").Value = 0.09898315953327597
$ws.Range("D8").Value = 0.1952872273609074
$ws.Range("E8").Value = 0.1671773588787815
$ws.Range("F8").Value = 1.386240639117979
$ws.Range("G8").Value = 0.002457589478532831
$ws.Range("I8").Value = 0.4976718254593342
$ws.Range("J8").Value = 0.1853462995587023
$ws.Range("K8").Value = 0.9059464436386975
$ws.Range("O8").Value = 3.314108322015414

# Row 9
$ws.Range("B9").Value = 0.1301302143922101
$ws.Range("D9").Value = 0.2106986858254345
$ws.Range("E9").Value = 0.1782226722325291
$ws.Range("F9").Value = 1.383656723864902
$ws.Range("G9").Value = 0.002450908866523327
$ws.Range("I9").Value = 0.4780690141069659
$ws.Range("J9").Value = 0.1959110609754191
$ws.Range("K9").Value = 1.170848122081964
$ws.Range("O9").Value = 3.272876426614346

# Row 10
$ws.Range("B10").Value = 0.1529262789605923
$ws.Range("D10").Value = 0.2223475566284208
$ws.Range("E10").Value = 0.1867305963239119
$ws.Range("F10").Value = 1.385989575911665
$ws.Range("G10").Value = 0.00244645112793237
$ws.Range("I10").Value = 0.4651729792018022
$ws.Range("J10").Value = 0.2042043063437404
$ws.Range("K10").Value = 1.364855375813136
$ws.Range("O10").Value = 3.253357573427166

# Row 11
$ws.Range("B11").Value = 0.1632759009097953
$ws.Range("D11").Value = 0.2277170569792872
$ws.Range("E11").Value = 0.1906864024394537
$ws.Range("F11").Value = 1.387972608482031
$ws.Range("G11").Value = 0.002444520032711406
$ws.Range("I11").Value = 0.4596333582241741
$ws.Range("J11").Value = 0.2080930827519722
$ws.Range("K11").Value = 1.452970513048911
$ws.Range("O11").Value = 3.246827415689864

# Row 12
$ws.Range("B12").Value = 0.1671918931639027
$ws.Range("D12").Value = 0.2297603805913866
$ws.Range("E12").Value = 0.1921966365587053
$ws.Range("F12").Value = 1.38885628577151
$ws.Range("G12").Value = 0.002443802616344461
$ws.Range("I12").Value = 0.4575826713856017
$ws.Range("J12").Value = 0.2095823861085364
$ws.Range("K12").Value = 1.486316069248517
$ws.Range("O12").Value = 3.244693171502888

# Row 13
$ws.Range("B13").Value = 0.1663486597531119
$ws.Range("D13").Value = 0.229319870808169
$ws.Range("E13").Value = 0.1918708361562551
$ws.Range("F13").Value = 1.388660063104339
$ws.Range("G13").Value = 0.002443956509798312
$ws.Range("I13").Value = 0.4580222302896679
$ws.Range("J13").Value = 0.2092608947634602
$ws.Range("K13").Value = 1.47913549971031
$ws.Range("O13").Value = 3.245137747917454

# Row 14
$ws.Range("B14").Value = 0.1635981371963027
$ws.Range("D14").Value = 0.2278849624426016
$ws.Range("E14").Value = 0.1908104049050721
$ws.Range("F14").Value = 1.388042647624857
$ws.Range("G14").Value = 0.002444460733346012
$ws.Range("I14").Value = 0.4594637038132907
$ws.Range("J14").Value = 0.2082152737358882
$ws.Range("K14").Value = 1.455714316961576
$ws.Range("O14").Value = 3.246645039858038

# Row 15
$ws.Range("B15").Value = 0.1619129405952151
$ws.Range("D15").Value = 0.227007340057142
$ws.Range("E15").Value = 0.1901624550805749
$ws.Range("F15").Value = 1.387681755920212
$ws.Range("G15").Value = 0.002444771386225601
$ws.Range("I15").Value = 0.4603527763292528
$ws.Range("J15").Value = 0.2075769767837272
$ws.Range("K15").Value = 1.441365287369933
$ws.Range("O15").Value = 3.247612416926302

# Row 16
$ws.Range("B16").Value = 0.152249474617534
$ws.Range("D16").Value = 0.2219980538773001
$ws.Range("E16").Value = 0.1864737923376012
$ws.Range("F16").Value = 1.385878545527433
$ws.Range("G16").Value = 0.002446579270944786
$ws.Range("I16").Value = 0.4655415874991347
$ws.Range("J16").Value = 0.2039525033025456
$ws.Range("K16").Value = 1.359093892072792
$ws.Range("O16").Value = 3.253831671523415

# Row 17
$ws.Range("B17").Value = 0.1463158434384155
$ws.Range("D17").Value = 0.2189429650980657
$ws.Range("E17").Value = 0.1842327896908529
$ws.Range("F17").Value = 1.385008567286448
$ws.Range("G17").Value = 0.002447713083485722
$ws.Range("I17").Value = 0.4688085022507522
$ws.Range("J17").Value = 0.2017587573167248
$ws.Range("K17").Value = 1.308586095676617
$ws.Range("O17").Value = 3.25824919871016

# Row 18
$ws.Range("B18").Value = 0.1429010682242478
$ws.Range("D18").Value = 0.2171923889371214
$ws.Range("E18").Value = 0.182951875113325
$ws.Range("F18").Value = 1.384594932178473
$ws.Range("G18").Value = 0.002448374333138651
$ws.Range("I18").Value = 0.4707183070325334
$ws.Range("J18").Value = 0.2005079036056543
$ws.Range("K18").Value = 1.279522291826368
$ws.Range("O18").Value = 3.261011101425169

# Row 19
$ws.Range("B19").Value = 0.1417445639894765
$ws.Range("D19").Value = 0.2166008159754824
$ws.Range("E19").Value = 0.1825195634064372
$ws.Range("F19").Value = 1.384469776584695
$ws.Range("G19").Value = 0.002448599787606254
$ws.Range("I19").Value = 0.471370216077232
$ws.Range("J19").Value = 0.2000862631977611
$ws.Range("K19").Value = 1.269679595135528
$ws.Range("O19").Value = 3.261984177014028

# Row 20
$ws.Range("B20").Value = 0.1469476878291402
$ws.Range("D20").Value = 0.2192674991135988
$ws.Range("E20").Value = 0.1844705153026922
$ws.Range("F20").Value = 1.38509219791834
$ws.Range("G20").Value = 0.00244759144459892
$ws.Range("I20").Value = 0.4684575495283729
$ws.Range("J20").Value = 0.2019911539535286
$ws.Range("K20").Value = 1.313964099912027
$ws.Range("O20").Value = 3.257756060987447

# Row 21
$ws.Range("B21").Value = 0.1644061207548759
$ws.Range("D21").Value = 0.2283061589992457
$ws.Range("E21").Value = 0.1911215469565875
$ws.Range("F21").Value = 1.388220393309354
$ws.Range("G21").Value = 0.002444312255229928
$ws.Range("I21").Value = 0.459039030868249
$ws.Range("J21").Value = 0.2085219444631008
$ws.Range("K21").Value = 1.462594288940579
$ws.Range("O21").Value = 3.246193116273531

# Row 22
$ws.Range("B22").Value = 0.175797537487739
$ws.Range("D22").Value = 0.2342717435303285
$ws.Range("E22").Value = 0.1955398006483691
$ws.Range("F22").Value = 1.391038683034935
$ws.Range("G22").Value = 0.002442249798705157
$ws.Range("I22").Value = 0.4531577349354485
$ws.Range("J22").Value = 0.2128875848769241
$ws.Range("K22").Value = 1.559605274294199
$ws.Range("O22").Value = 3.240609962136062

# Row 23
$ws.Range("B23").Value = 0.1697195146651325
$ws.Range("D23").Value = 0.2310825004140469
$ws.Range("E23").Value = 0.1931751737218121
$ws.Range("F23").Value = 1.389463639868069
$ws.Range("G23").Value = 0.002443343209789006
$ws.Range("I23").Value = 0.4562715836798175
$ws.Range("J23").Value = 0.210548646590027
$ws.Range("K23").Value = 1.507840885025644
$ws.Range("O23").Value = 3.243408923190657

# Row 24
$ws.Range("B24").Value = 0.1466620417126876
$ws.Range("D24").Value = 0.2191207591027933
$ws.Range("E24").Value = 0.1843630162906109
$ws.Range("F24").Value = 1.38505411899466
$ws.Range("G24").Value = 0.002447646408243709
$ws.Range("I24").Value = 0.4686161167884908
$ws.Range("J24").Value = 0.2018860551655166
$ws.Range("K24").Value = 1.311532785681777
$ws.Range("O24").Value = 3.257978316436407

# Row 25
$ws.Range("B25").Value = 0.1217188263515112
$ws.Range("D25").Value = 0.2064719574093061
$ws.Range("E25").Value = 0.1751655991268208
$ws.Range("F25").Value = 1.383613633125336
$ws.Range("G25").Value = 0.00245263670537474
$ws.Range("I25").Value = 0.4831076251381168
$ws.Range("J25").Value = 0.1929598852404979
$ws.Range("K25").Value = 1.09928955296806
$ws.Range("O25").Value = 3.282141937389525
